$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three now-obsolete trailing rows (18-FEB-26 / 21-FEB-26 groupings no longer present)
$ws.Range("A17:K19").EntireRow.Delete()

# Keep the Date column stored as text (matches source data authored as strings, not Excel dates)
$ws.Range("A2:A16").NumberFormat = "@"

$reportData = @(
  @{ Row = 2; A = '31-JAN-26'; B = 'SM-436'; C = 'Air Arabia Egypt E5-592'; D = 344; E = 663; F = -319; G = 30; H = 30; I = 0; J = 'LOW THREAT'; K = 'SAR' }
  @{ Row = 3; A = '31-JAN-26'; B = 'SM-436'; C = 'Nile Air NP-116'; D = 345; E = 663; F = -318; G = 30; H = 30; I = 0; J = 'LOW THREAT'; K = 'SAR' }
  @{ Row = 4; A = '31-JAN-26'; B = 'SM-436'; C = 'EgyptAir MS-634'; D = 524; E = 663; F = -139; G = 46; H = 30; I = -16; J = 'LOW THREAT'; K = 'SAR' }
  @{ Row = 5; A = '04-FEB-26'; B = 'SM-436'; C = 'Air Arabia Egypt E5-592'; D = 345; E = 663; F = -318; G = 30; H = 30; I = 0; J = 'LOW THREAT'; K = 'SAR' }
  @{ Row = 6; A = '04-FEB-26'; B = 'SM-436'; C = 'Nile Air NP-106'; D = 401; E = 663; F = -262; G = 30; H = 30; I = 0; J = 'LOW THREAT'; K = 'SAR' }
  @{ Row = 7; A = '07-FEB-26'; B = 'SM-436'; C = 'Nile Air NP-116'; D = 350; E = 663; F = -313; G = 30; H = 30; I = 0; J = 'LOW THREAT'; K = 'SAR' }
  @{ Row = 8; A = '07-FEB-26'; B = 'SM-436'; C = 'Air Arabia Egypt E5-592'; D = 370; E = 663; F = -293; G = 30; H = 30; I = 0; J = 'LOW THREAT'; K = 'SAR' }
  @{ Row = 9; A = '07-FEB-26'; B = 'SM-436'; C = 'Nesma Airlines NE-141'; D = 400; E = 663; F = -263; G = 30; H = 30; I = 0; J = 'LOW THREAT'; K = 'SAR' }
  @{ Row = 10; A = '14-FEB-26'; B = 'SM-436'; C = 'Nesma Airlines NE-141'; D = 400; E = 895; F = -495; G = 30; H = 30; I = 0; J = 'MEDIUM THREAT - MONITOR'; K = 'SAR' }
  @{ Row = 11; A = '14-FEB-26'; B = 'SM-436'; C = 'Air Arabia Egypt E5-592'; D = 407; E = 895; F = -488; G = 30; H = 30; I = 0; J = 'MEDIUM THREAT - MONITOR'; K = 'SAR' }
  @{ Row = 12; A = '14-FEB-26'; B = 'SM-436'; C = 'Nile Air NP-116'; D = 521; E = 895; F = -374; G = 30; H = 30; I = 0; J = 'LOW THREAT'; K = 'SAR' }
  @{ Row = 13; A = '14-FEB-26'; B = 'SM-436'; C = 'EgyptAir MS-634'; D = 732; E = 895; F = -163; G = 46; H = 30; I = -16; J = 'MEDIUM THREAT - MONITOR'; K = 'SAR' }
  @{ Row = 14; A = '25-FEB-26'; B = 'SM-436'; C = 'Nesma Airlines NE-141'; D = 350; E = 669; F = -319; G = 30; H = 30; I = 0; J = 'LOW THREAT'; K = 'SAR' }
  @{ Row = 15; A = '25-FEB-26'; B = 'SM-436'; C = 'Nile Air NP-106'; D = 350; E = 669; F = -319; G = 30; H = 30; I = 0; J = 'LOW THREAT'; K = 'SAR' }
  @{ Row = 16; A = '25-FEB-26'; B = 'SM-436'; C = 'Air Arabia Egypt E5-592'; D = 364; E = 669; F = -305; G = 30; H = 30; I = 0; J = 'LOW THREAT'; K = 'SAR' }
)

foreach ($row in $reportData) {
  $r = $row.Row
  $ws.Cells.Item($r, 1).Value = $row.A
  $ws.Cells.Item($r, 2).Value = $row.B
  $ws.Cells.Item($r, 3).Value = $row.C
  $ws.Cells.Item($r, 4).Value = $row.D
  $ws.Cells.Item($r, 5).Value = $row.E
  $ws.Cells.Item($r, 6).Value = $row.F
  $ws.Cells.Item($r, 7).Value = $row.G
  $ws.Cells.Item($r, 8).Value = $row.H
  $ws.Cells.Item($r, 9).Value = $row.I
  $ws.Cells.Item($r, 11).Value = $row.K
}

# Recolor the IMPACT column (J) to match LOW THREAT (green) / MEDIUM THREAT (yellow) styling;
# the HIGH THREAT (red) look is no longer used anywhere in the refreshed report
$lowTemplate = $ws.Range("J2")
$mediumTemplate = $ws.Range("J4")
foreach ($row in $reportData) {
  $r = $row.Row
  $target = $ws.Cells.Item($r, 10)
  if ($row.J -eq "MEDIUM THREAT - MONITOR") {
    $mediumTemplate.Copy()
  } else {
    $lowTemplate.Copy()
  }
  $target.PasteSpecial(-4122)
  $target.Value = $row.J
}
$excel.CutCopyMode = $false

# Column J no longer needs to fit the long HIGH THREAT label, narrow it back down
$ws.Columns.Item(10).ColumnWidth = 25

$ws.Range("A1").Select()
